$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: the "FinSimboloIncPos" labels become "FinSimbolo", except I8 ---
# (I8 keeps its original "FinSimboloIncPos" value; every other cell B8:R8 changes)
$ws.Range("B8:H8").Value2 = "FinSimbolo"
$ws.Range("J8:N8").Value2 = "FinSimbolo"
$ws.Range("O8:R8").Value2 = "FinSimbolo"

# O8 had picked up a slightly different (but visually identical) cell format than
# its row neighbours; copy N8's format onto it so it is consistent with the rest
# of the row, same as the others in row 8.
$ws.Range("N8").Copy()
$ws.Range("O8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- View state: zoom out to 62% and move the selection to the row we just edited ---
$excel.ActiveWindow.Zoom = 62
$ws.Range("B8:R8").Select()
